$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  13"
$ws.Range("C9").Value = "Report Covering the Week  3/27/2023  Through  4/2/2023"

# --- Cells that flip between the "0"/"***.*" text placeholders and real numbers ---
# (Copy() from an existing cell of the desired style/type first, so the style index
#  lands on the same shared style as similar cells elsewhere in the table; then overwrite
#  the value for the cells that need to hold a real number rather than the placeholder text.)
$ws.Range("F19").Copy($ws.Range("C18"))
$ws.Range("C18").Value = 3

$ws.Range("C14").Copy($ws.Range("C22"))

$ws.Range("F19").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 2

$ws.Range("M22").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100

$ws.Range("F19").Copy($ws.Range("G22"))
$ws.Range("G22").Value = 2

$ws.Range("M22").Copy($ws.Range("H22"))
$ws.Range("H22").Value = -50

$ws.Range("C14").Copy($ws.Range("D26"))

$ws.Range("E14").Copy($ws.Range("E26"))

# --- Remaining weekly crime-stat figures (straightforward numeric overwrites) ---
# Row 14
$ws.Range("F14").Value = 1
# Row 16
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = 21.428571428571
$ws.Range("I16").Value = 35
$ws.Range("J16").Value = 52
$ws.Range("K16").Value = -32.692307692307
$ws.Range("L16").Value = 66.666666666666
$ws.Range("M16").Value = -51.388888888888
$ws.Range("N16").Value = -88.709677419354
# Row 17
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 32
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 95
$ws.Range("J17").Value = 77
$ws.Range("K17").Value = 23.376623376623
$ws.Range("L17").Value = 69.642857142857
$ws.Range("M17").Value = 82.692307692307
$ws.Range("N17").Value = 23.376623376623
# Row 18
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -41.666666666666
$ws.Range("I18").Value = 26
$ws.Range("J18").Value = 29
$ws.Range("K18").Value = -10.344827586206
$ws.Range("L18").Value = -16.129032258064
$ws.Range("M18").Value = -61.764705882352
$ws.Range("N18").Value = -92.307692307692
# Row 19
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -14.285714285714
$ws.Range("G19").Value = 25
$ws.Range("H19").Value = -36
$ws.Range("I19").Value = 78
$ws.Range("J19").Value = 102
$ws.Range("K19").Value = -23.529411764705
$ws.Range("L19").Value = 50
$ws.Range("M19").Value = 2.631578947368
$ws.Range("N19").Value = -39.534883720930
# Row 20
$ws.Range("C20").Value = 3
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 67
$ws.Range("J20").Value = 60
$ws.Range("K20").Value = 11.666666666666
$ws.Range("L20").Value = 86.111111111111
$ws.Range("M20").Value = 8.064516129032
$ws.Range("N20").Value = -92.756756756756
# Row 21
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = 9.090909090909
$ws.Range("F21").Value = 88
$ws.Range("G21").Value = 93
$ws.Range("H21").Value = -5.376344086021
$ws.Range("I21").Value = 305
$ws.Range("J21").Value = 323
$ws.Range("K21").Value = -5.572755417956
$ws.Range("L21").Value = 48.780487804878
$ws.Range("M21").Value = -7.854984894259
$ws.Range("N21").Value = -82.951369480156
# Row 22
$ws.Range("J22").Value = 7
$ws.Range("K22").Value = -42.857142857142
# Row 24
$ws.Range("C24").Value = 35
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = 34.615384615384
$ws.Range("G24").Value = 104
$ws.Range("H24").Value = 16.346153846153
$ws.Range("I24").Value = 382
$ws.Range("J24").Value = 333
$ws.Range("K24").Value = 14.714714714714
$ws.Range("L24").Value = 35.943060498220
$ws.Range("M24").Value = 88.177339901477
# Row 25
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = -33.333333333333
$ws.Range("G25").Value = 44
$ws.Range("H25").Value = 13.636363636363
$ws.Range("I25").Value = 139
$ws.Range("J25").Value = 121
$ws.Range("K25").Value = 14.876033057851
$ws.Range("L25").Value = 52.747252747252
$ws.Range("M25").Value = -3.472222222222
# Row 27
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 20
$ws.Range("I27").Value = 16
$ws.Range("J27").Value = 13
$ws.Range("K27").Value = 23.076923076923
$ws.Range("L27").Value = 60
# Row 28
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -25
$ws.Range("L28").Value = -20
# Row 29
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 50
$ws.Range("L29").Value = 33.333333333333
